# Adds a UUID field to each form/section header, and introduces a new
# "Jurisdiction" / JURISDICTION field on the two rows that reference a
# jurisdiction (Jurisdiction-Agency and Jurisdiction-Forest Monitoring
# Datasets sections).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Structural change: insert the two new "Jurisdiction" rows ------------
# Row 30 (old "Jurisdiction-Agency" header) shifts down to make room for a
# new C/D pair that records the jurisdiction name's field type.
$ws.Rows.Item(30).Insert()
# After the first insert, the old "Name of Jurisdictional Agency" row has
# moved to row 33; insert again before it for the second new field row.
$ws.Rows.Item(33).Insert()

$ws.Cells.Item(30, 3).Value = "Jurisdiction"
$ws.Cells.Item(30, 4).Value = "JURISDICTION"
$ws.Cells.Item(33, 3).Value = "Jurisdiction"
$ws.Cells.Item(33, 4).Value = "JURISDICTION"

# --- New column F: a UUID per form/section header --------------------------
$ws.Cells.Item(2, 6).Value  = "d6493bb0-a610-11e6-a350-8d1f5c467d3b"
$ws.Cells.Item(8, 6).Value  = "d6493bb1-a610-11e6-a350-8d1f5c467d3b"
$ws.Cells.Item(13, 6).Value = "d6493bb2-a610-11e6-a350-8d1f5c467d3b"
$ws.Cells.Item(18, 6).Value = "d6493bb3-a610-11e6-a350-8d1f5c467d3b"
$ws.Cells.Item(23, 6).Value = "d6493bb4-a610-11e6-a350-8d1f5c467d3b"
$ws.Cells.Item(28, 6).Value = "d6493bb5-a610-11e6-a350-8d1f5c467d3b"
$ws.Cells.Item(31, 6).Value = "d6493bb6-a610-11e6-a350-8d1f5c467d3b"
$ws.Cells.Item(36, 6).Value = "d6493bb7-a610-11e6-a350-8d1f5c467d3b"

# --- Cosmetic: column widths to fit the two new (wider) columns -----------
$ws.Columns.Item(1).ColumnWidth = 7.671768707482998
$ws.Columns.Item(2).ColumnWidth = 7.671768707482998
$ws.Columns.Item(3).ColumnWidth = 28.457482993197267
$ws.Columns.Item(4).ColumnWidth = 13.069727891156466
$ws.Columns.Item(5).ColumnWidth = 63.258503401360564
$ws.Columns.Item(6).ColumnWidth = 34.31972789115646
$ws.Columns.Item(7).ColumnWidth = 53.48809523809526

# --- Cosmetic: move the active selection to G1, as in the target file -----
$ws.Range("G1").Select() | Out-Null
